$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1E - Constant")

# Update G3 formula (trial 3, row for "Learning Rate"/epoch-time style calc)
$ws.Range("G3").Formula = "=3052.63/60"

# Fill in G4:G11 with new trial values
$ws.Range("G4").Value = 7247.31
$ws.Range("G5").Value = 6808.2
$ws.Range("G6").Value = 6220.42
$ws.Range("G7").Value = 6007.96
$ws.Range("G8").Value = 5893.071
$ws.Range("G9").Value = 5816.88
$ws.Range("G10").Value = 5760.135
$ws.Range("G11").Value = 5722.85

# Update G12 value
$ws.Range("G12").Value = 5677.75634

# Update the active selection to G3
$ws.Range("G3").Select()
